# Automatische test-sync: 2025-08-06 20:01:50
# Appends the new mail-log entry (row 12) to the "Logs" sheet, extends the
# conditional formatting ranges that covered rows 2-11 so they also cover
# row 12, and bumps the "Inkoop / Bestellingen" tally on the "Dashboard"
# sheet from 4 to 5.

$wb = $excel.ActiveWorkbook
$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

# --- Add the new log row -----------------------------------------------
$newRow = 12
$logs.Range("A" + $newRow).Value = "Is dit artikel momenteel beschikbaar?"
$logs.Range("B" + $newRow).Value = "mailmind.test@zohomail.eu"
$logs.Range("C" + $newRow).Value = "Testmail #1: Is dit artikel momenteel beschikbaar?"
$logs.Range("D" + $newRow).Value = "Inkoop / Bestellingen"
$logs.Range("E" + $newRow).Value = "Bedankt, we hebben dit doorgestuurd naar inkoop@bedrijf.nl."
$logs.Range("F" + $newRow).Value = "2025-08-06 20:01:12"
$logs.Range("G" + $newRow).Value = "Ja"
$logs.Range("H" + $newRow).Value = "Ja"
$logs.Range("I" + $newRow).Value = "Nee"
$logs.Range("J" + $newRow).Value = "Nee"

# --- Extend conditional formatting ranges to include the new row -------
$ccols = @("D", "G", "H", "I", "J")
foreach ($col in $ccols) {
    $oldRange = $logs.Range($col + "2:" + $col + "11")
    $newRange = $logs.Range($col + "2:" + $col + $newRow)
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}

# --- Update the Dashboard roll-up count ---------------------------------
$dashboard.Range("B3").Value = 5
